$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Price/Volume columns to Text so numeric-looking
# strings (e.g. "0.9999", "13.14") are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.195.83'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '1.858.46'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '0.7131'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').Value = '240.32'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '0.07756'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('D9').Value = '0.3072'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '25.09'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').Value = '0.08253'
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.849.59'
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '5.228'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.7155'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').Value = '29.170.46'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').Value = '5.855'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '243.75'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').Value = '0.000007785'
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('D20').Value = '13.14'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').Value = '2.104.01'
$ws.Range('E21').Value = '  -1.54%  '
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '7.972'
$ws.Range('E23').Value = '  +2.41%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = '0.1595'
$ws.Range('E25').Value = '  +2.86%  '
$ws.Range('D26').Value = '162.42'
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').Value = '8.900'
$ws.Range('E27').Value = '  -1.20%  '
$ws.Range('D28').Value = '18.30'
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('D30').Value = '1.315'
$ws.Range('E30').Value = '  -3.08%  '
$ws.Range('D31').Value = '4.400'
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').Value = '4.210'
$ws.Range('E32').Value = '  +2.86%  '
$ws.Range('D33').Value = '0.05183'
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').Value = '1.909'
$ws.Range('E34').Value = '  -1.24%  '
$ws.Range('D35').Value = '1.170'
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D36').Value = '0.7261'
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').Value = '0.01852'
$ws.Range('E38').Value = '  -0.83%  '
$ws.Range('D39').Value = '2.685'
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('D40').Value = '1.165.13'
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('D41').Value = '0.9035'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('D42').Value = '6.155'
$ws.Range('E42').Value = '  +2.26%  '
$ws.Range('D43').Value = '72.16'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').Value = '0.9997'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').Value = '101.56'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('D46').Value = '2.001.21'
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('D47').Value = '0.5213'
$ws.Range('E47').Value = '  -2.82%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  -1.23%  '
$ws.Range('D50').Value = '9.307'
$ws.Range('E50').Value = '  +1.71%  '
$ws.Range('D51').Value = '2.856'
$ws.Range('E51').Value = '  +1.21%  '

# Restore original (default/general) formatting so no stray styles remain.
$ws.Range("D2:E51").ClearFormats()
